$wb = $excel.ActiveWorkbook

# --- Update status / timestamp values to reflect "Ready for handoff" ---

# Overview sheet: E2/F2 = Status columns for zh-cn / de-de, G2 = Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-18 16:58:00"

# zh-cn sheet: C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-18 16:57:54"

# de-de sheet: C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-18 16:58:00"

# --- Shrink the Status columns now that the text is shorter ---
# (ColumnWidth is internally quantized to 1/6-character steps by the engine,
#  so the input below is chosen to land on the nearest achievable step to 17.22)
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
